$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear rows 10-29 first so the sheet dimension shrinks to A1:A9
$ws.Range("A10:A29").Clear()

# Consolidate token text fields into tuple-string representation
$ws.Range("A2").Value = '(''Angel'', [''Token Creature — Angel'', ''Flying, vigilance'', ''4/4''])'
$ws.Range("A3").Value = '(''Bird Illusion'', [''Token Creature — Bird Illusion'', ''Flying'', ''1/1''])'
$ws.Range("A4").Value = '(''Elf Knight'', [''Token Creature — Elf Knight'', ''Vigilance'', ''2/2''])'
$ws.Range("A5").Value = '(''Goblin'', [''Token Creature — Goblin'', ''1/1''])'
$ws.Range("A6").Value = '(''Insect'', [''Token Creature — Insect'', ''1/1''])'
$ws.Range("A7").Value = '(''Ral, Izzet Viceroy Emblem'', [''Emblem — Ral'', ''Whenever you cast an instant or sorcery spell, this emblem deals 4 damage to any target and you draw two cards.''])'
$ws.Range("A8").Value = '(''Soldier'', [''Token Creature — Soldier'', ''Lifelink'', ''1/1''])'
$ws.Range("A9").Value = '(''Vraska, Golgari Queen Emblem'', [''Emblem — Vraska'', ''Whenever a creature you control deals combat damage to a player, that player loses the game.''])'
